$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text: "In Translation" -> "Ready for handoff" (every Status cell)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest Handoff Date/Datetime timestamps (stored as text, keep as text)
$overview.Range("D2").Value = "2016-03-21 10:32:34"
$dede.Range("E2").Value = "2016-03-21 10:32:34"
$zhcn.Range("E2").Value = "2016-03-21 10:32:31"
